# Convert {{placeholder}} style merge fields to [%placeholder%] style
# used by the new templating engine.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "{{dissertation_topic}}"; New = "[%dissertation_topic%]" },
    @{ Old = "{{student_specialty}}"; New = "[%student_specialty%]" },
    @{ Old = "{{student_full_name}}"; New = "[%student_full_name%]" },
    @{ Old = "{{student_program}}"; New = "[%student_program%]" },
    @{ Old = "{{student_phone}}"; New = "[%student_phone%]" },
    @{ Old = "{{student_email}}"; New = "[%student_email%]" },
    @{ Old = "{{day}}"; New = "[%day%]" },
    @{ Old = "{{month}}"; New = "[%month%]" },
    @{ Old = "{{year}}"; New = "[%year%]" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.New, 2)
}
